$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 113.5
$ws.Range("I33").Value = 117.75
$ws.Range("K33").Value = 117.75
$ws.Range("M33").Value = 111.25

$ws.Range("H40").Value = 1666.6666
$ws.Range("I40").Value = 1328.5714
$ws.Range("J40").Value = 2140
$ws.Range("K40").Value = 1328.5714
$ws.Range("L40").Value = 2140
$ws.Range("M40").Value = -1153.5714
$ws.Range("N40").Value = -2490

$ws.Range("H98").Value = 490.4091
$ws.Range("I98").Value = 425.42105
$ws.Range("J98").Value = 902
$ws.Range("K98").Value = 425.42105
$ws.Range("L98").Value = 902
$ws.Range("M98").Value = 1072.57895
$ws.Range("N98").Value = -3898

$ws.Range("H106").Value = 12822592
$ws.Range("I106").Value = 18519886
$ws.Range("K106").Value = 18519886
$ws.Range("M106").Value = -18519255

$ws.Range("H113").Value = 50003690
$ws.Range("J113").Value = 4891.2856
$ws.Range("L113").Value = 4891.2856
$ws.Range("N113").Value = -11399.2856

$ws.Range("H116").Value = 3434.25
$ws.Range("I116").Value = 1377.4
$ws.Range("K116").Value = 1377.4
$ws.Range("M116").Value = 2064.6

$ws.Range("H122").Value = 490.4091
$ws.Range("I122").Value = 425.42105
$ws.Range("J122").Value = 902
$ws.Range("K122").Value = 1276.26315
$ws.Range("L122").Value = 2706
$ws.Range("M122").Value = 1173.73685
$ws.Range("N122").Value = -7606

$ws.Range("H125").Value = 980
$ws.Range("I125").Value = 980
$ws.Range("K125").Value = 8820
$ws.Range("M125").Value = -6360

$ws.Range("H129").Value = 213573.52
$ws.Range("J129").Value = 238958.53
$ws.Range("L129").Value = 716875.59
$ws.Range("N129").Value = -726875.59

$ws.Range("H132").Value = 2000
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()

$ws.Range("H137").Value = 2767.0952
$ws.Range("I137").Value = 2665.2354
$ws.Range("K137").Value = 7995.706200000001
$ws.Range("M137").Value = -5445.706200000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3750.3635
$ws.Range("I32").Value = 3821.4753
$ws.Range("J32").Value = 2882.8
$ws.Range("K32").Value = 3821.4753
$ws.Range("L32").Value = 2882.8
$ws.Range("M32").Value = -3534.4753
$ws.Range("N32").Value = -3456.8

$ws.Range("H61").Value = 2305.3044
$ws.Range("I61").Value = 1953.4286
$ws.Range("J61").Value = 6000
$ws.Range("K61").Value = 1953.4286
$ws.Range("L61").Value = 6000
$ws.Range("M61").Value = -1741.4286
$ws.Range("N61").Value = -6424

$ws.Range("H74").Value = 166667860
$ws.Range("I74").Value = 200001020
$ws.Range("J74").Value = 2000
$ws.Range("K74").Value = 200001020
$ws.Range("L74").Value = 2000
$ws.Range("M74").Value = -200000146
$ws.Range("N74").Value = -3748

$ws.Range("H77").Value = 166667860
$ws.Range("I77").Value = 200001020
$ws.Range("J77").Value = 2000
$ws.Range("K77").Value = 1000005100
$ws.Range("L77").Value = 10000
$ws.Range("M77").Value = -1000000732
$ws.Range("N77").Value = -18736

$ws.Range("H132").Value = 13078.75
$ws.Range("I132").Value = 1551.0834
$ws.Range("J132").Value = 64953.25
$ws.Range("K132").Value = 4653.2502
$ws.Range("L132").Value = 194859.75
$ws.Range("M132").Value = -2123.2502
$ws.Range("N132").Value = -199919.75

$ws.Range("H136").Value = 2305.3044
$ws.Range("I136").Value = 1953.4286
$ws.Range("J136").Value = 6000
$ws.Range("K136").Value = 5860.2858
$ws.Range("L136").Value = 18000
$ws.Range("M136").Value = -3310.2858
$ws.Range("N136").Value = -23100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1725.5625
$ws.Range("I99").Value = 1159.8
$ws.Range("J99").Value = 2668.5
$ws.Range("K99").Value = 1159.8
$ws.Range("L99").Value = 2668.5
$ws.Range("M99").Value = 338.2
$ws.Range("N99").Value = -5664.5

$ws.Range("H105").Value = 2348.353
$ws.Range("I105").Value = 1591.3334
$ws.Range("J105").Value = 3200
$ws.Range("K105").Value = 1591.3334
$ws.Range("L105").Value = 3200
$ws.Range("M105").Value = 155.6666
$ws.Range("N105").Value = -6694

$ws.Range("H134").Value = 3707.6128
$ws.Range("I134").Value = 3859.862
$ws.Range("J134").Value = 1500
$ws.Range("K134").Value = 11579.586
$ws.Range("L134").Value = 4500
$ws.Range("M134").Value = -9044.585999999999
$ws.Range("N134").Value = -9570

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 423.7143
$ws.Range("I22").Value = 304.2857
$ws.Range("J22").Value = 543.1429000000001
$ws.Range("K22").Value = 304.2857
$ws.Range("L22").Value = 543.1429000000001
$ws.Range("M22").Value = 45.71429999999998
$ws.Range("N22").Value = -1243.1429

$ws.Range("H31").Value = 14688.429
$ws.Range("I31").Value = 26637.562
$ws.Range("J31").Value = 4626
$ws.Range("K31").Value = 26637.562
$ws.Range("L31").Value = 4626
$ws.Range("M31").Value = -26342.562
$ws.Range("N31").Value = -5216

$ws.Range("H34").Value = 14688.429
$ws.Range("I34").Value = 26637.562
$ws.Range("J34").Value = 4626
$ws.Range("K34").Value = 26637.562
$ws.Range("L34").Value = 4626
$ws.Range("M34").Value = -26435.562
$ws.Range("N34").Value = -5030

$ws.Range("H58").Value = 14925.75
$ws.Range("I58").Value = 1052.88
$ws.Range("J58").Value = 46455
$ws.Range("K58").Value = 1052.88
$ws.Range("L58").Value = 46455
$ws.Range("M58").Value = -849.8800000000001
$ws.Range("N58").Value = -46861

$ws.Range("H132").Value = 11091.8545
$ws.Range("I132").Value = 13959.075
$ws.Range("K132").Value = 41877.22500000001
$ws.Range("M132").Value = -39347.22500000001

$ws.Range("H134").Value = 1034.5074
$ws.Range("I134").Value = 763.5909
$ws.Range("K134").Value = 2290.7727
$ws.Range("M134").Value = 244.2273

$ws.Range("H136").Value = 14925.75
$ws.Range("I136").Value = 1052.88
$ws.Range("J136").Value = 46455
$ws.Range("K136").Value = 3158.64
$ws.Range("L136").Value = 139365
$ws.Range("M136").Value = -608.6400000000003
$ws.Range("N136").Value = -144465

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H42").Value = 2000
$ws.Range("J42").Value = 2000
$ws.Range("L42").Value = 6000
$ws.Range("N42").Value = -7068

$ws.Range("H81").Value = 4843.3335
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 4843.3335
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 14530.0005
$ws.Range("M81").ClearContents()
$ws.Range("N81").Value = -16776.0005

$ws.Range("H84").Value = 4843.3335
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 4843.3335
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 43590.0015
$ws.Range("M84").ClearContents()
$ws.Range("N84").Value = -54822.0015

$ws.Range("H113").Value = 502.32
$ws.Range("I113").Value = 426.07693
$ws.Range("J113").Value = 584.9167
$ws.Range("K113").Value = 1278.23079
$ws.Range("L113").Value = 1754.7501
$ws.Range("M113").Value = 891.7692099999999
$ws.Range("N113").Value = -6094.7501

$ws.Range("H131").Value = 113157.49
$ws.Range("J131").Value = 115743.49
$ws.Range("L131").Value = 347230.47
$ws.Range("N131").Value = -357310.47

$ws.Range("H137").Value = 23812082
$ws.Range("I137").Value = 796.6667
$ws.Range("J137").Value = 41670548
$ws.Range("K137").Value = 2390.0001
$ws.Range("L137").Value = 125011644
$ws.Range("M137").Value = 2709.9999
$ws.Range("N137").Value = -125021844

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6667.778
$ws.Range("I70").Value = 7600
$ws.Range("J70").Value = 5502.5
$ws.Range("K70").Value = 7600
$ws.Range("L70").Value = 5502.5
$ws.Range("M70").Value = -7330
$ws.Range("N70").Value = -6042.5

$ws.Range("H73").Value = 6667.778
$ws.Range("I73").Value = 7600
$ws.Range("J73").Value = 5502.5
$ws.Range("K73").Value = 7600
$ws.Range("L73").Value = 5502.5
$ws.Range("M73").Value = -6664
$ws.Range("N73").Value = -7374.5

$ws.Range("H132").Value = 20994.104
$ws.Range("I132").Value = 3727.4348
$ws.Range("J132").Value = 87183
$ws.Range("K132").Value = 11182.3044
$ws.Range("L132").Value = 261549
$ws.Range("M132").Value = -8652.304400000001
$ws.Range("N132").Value = -266609

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 936132.4399999999
$ws.Range("I122").Value = 1309405.4
$ws.Range("K122").Value = 3928216.2
$ws.Range("M122").Value = -3925766.2

$ws.Range("H132").Value = 2196.2683
$ws.Range("I132").Value = 1091
$ws.Range("J132").Value = 3608.5557
$ws.Range("K132").Value = 3273
$ws.Range("L132").Value = 10825.6671
$ws.Range("M132").Value = -743
$ws.Range("N132").Value = -15885.6671

$ws.Range("H136").Value = 34413.4
$ws.Range("I136").Value = 39554.31
$ws.Range("J136").Value = 997.5
$ws.Range("K136").Value = 118662.93
$ws.Range("L136").Value = 2992.5
$ws.Range("M136").Value = -116112.93
$ws.Range("N136").Value = -8092.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 5668
$ws.Range("I5").Value = 3501
$ws.Range("J5").Value = 10002
$ws.Range("K5").Value = 3501
$ws.Range("L5").Value = 10002
$ws.Range("M5").Value = -3389
$ws.Range("N5").Value = -10226

$ws.Range("H18").Value = 28000
$ws.Range("J18").Value = 28000
$ws.Range("L18").Value = 28000
$ws.Range("N18").Value = -28346

$ws.Range("H92").Value = 25525.25
$ws.Range("J92").Value = 25525.25
$ws.Range("L92").Value = 25525.25
$ws.Range("N92").Value = -30517.25

$ws.Range("H100").Value = 300
$ws.Range("I100").Value = 300
$ws.Range("K100").Value = 600
$ws.Range("M100").Value = -59

$ws.Range("H132").Value = 1521.0476
$ws.Range("I132").Value = 1132.9333
$ws.Range("K132").Value = 3398.7999
$ws.Range("M132").Value = -868.7999

$ws.Range("H136").Value = 37039236
$ws.Range("I136").Value = 62502076
$ws.Range("K136").Value = 187506228
$ws.Range("M136").Value = -187503678
